# Port-level landings cleanup for Table36.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 currently holds the "MONTEREY AREA TOTALS" label in B2 with no
# port name in A2. Pull the label over to A2 (matching the port-name
# column used by every other block) and replace B2 with the "Totals"
# label used by the other area-total rows.
$ws.Range("A2").Value2 = $ws.Range("B2").Value2
$ws.Range("B2").Value2 = "Totals"

# Columns A and B should share the same (wider) best-fit width now that
# column A carries the longer port/area labels too.
$ws.Range("A:B").ColumnWidth = 23

# Reset the view: scroll back to the top and select B4 instead of the
# previous mid-sheet selection.
$ws.Range("B4").Select()
